$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.926.78"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "3.440.99"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "3.438.66"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "4.035.82"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.33%  "
$ws.Range("D16").Value = "65.887.81"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "3.441.61"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.00%  "
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "2.777.41"
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0684"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "326.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.17%  "
